# Fix minutes/seconds formatting in the "Общее время" (total time) column.
# Strings look like "18 ч. 25 мин. 4 сек." and should become
# "18 ч. 25 мин. 04 сек." -- i.e. zero-pad the minutes and seconds
# components to two digits, while leaving the hours component untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

$changed = 0
for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $s = $cell.Text

    if ($s -match '^(\d+) ч\. (\d+) мин\. (\d+) сек\.$') {
        $hours = $matches[1]
        $minutes = $matches[2].PadLeft(2, '0')
        $seconds = $matches[3].PadLeft(2, '0')
        $newValue = "$hours ч. $minutes мин. $seconds сек."

        if ($newValue -ne $s) {
            $cell.Value = $newValue
            $changed = $changed + 1
        }
    }
}

Write-Host "Updated" $changed "haul time cells"
